$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B:G (ME, MAE, MSE, RMSE, SE, N)
$data = @(
    @(0.1887967460050058, 2.743091716767341, 23.88233659257806, 4.886955759220464, 4.931898846174919, 51),
    @(0.352110747706186,  2.607870699580363, 22.56754235743245, 4.750530744814989, 4.785560805511427, 50),
    @(0.2174705470141381, 2.756943539436507, 23.51947414538977, 4.849688046193257, 4.895016289175286, 49),
    @(0.4431031871169503, 2.719424870216065, 23.17433994440393, 4.813973405037041, 4.844263939605876, 48),
    @(0.2052815252275707, 2.707768026847837, 24.13652786648203, 4.912894041853745, 4.961670927367957, 47),
    @(0.4590069508388199, 2.706948687152313, 24.1233625610406,  4.911553986371381, 4.944094250912167, 46),
    @(0.04636964274274196,2.422883991970613, 21.4973978904842,  4.636528646572153, 4.688685967794817, 45),
    @(0.2175654990191061, 2.393845526388351, 21.18684633999888, 4.60291715545684,  4.650927664705682, 44),
    @(0.02747142410952863,2.403862371068361, 21.54546590760865, 4.641709373453777, 4.696560501412658, 43),
    @(0.2558878751990361, 2.356132134174123, 21.26811481140563, 4.611736637255605, 4.660447797710927, 42)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
